$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "id_building_type" column (column B) entirely - shifts everything left.
$ws.Columns("B").Delete()

# The former "id_cooling_technology" column is now column B; rename it (header +
# the underlying table column name) to "id_ventilation_technology".
$ws.Range("B1").Value = "id_ventilation_technology"

# Update the data row: id_region changes from 1 to 9, and the
# id_ventilation_technology value changes from 11 to 1.
$ws.Range("A2").Value = 9
$ws.Range("B2").Value = 1

# Resync the table's bounds/name/style so xl/tables/table1.xml matches the new
# (one-column-narrower) range.
$lo = $ws.ListObjects.Item(1)
$loName = $lo.Name
$lo.Unlist()
$lo2 = $ws.ListObjects.Add(1, $ws.Range("A1:AR2"), 0, 1)
$lo2.Name = $loName
$lo2.TableStyle = "TableStyleMedium6"

# Match the reset selection seen after the edit.
$ws.Range("A3").Select()
